$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.406155467033386
$ws.Range("B1").Value = 3.235643625259399
$ws.Range("C1").Value = 5.827722549438477
$ws.Range("D1").Value = 4.910142421722412
$ws.Range("E1").Value = 1.179466128349304
